# Apply the edit described by the diff:
#  - rename the shared string used by L1 from "Categoria_do_produto_pesquisado"
#    to "CATEGORIA_PRODUTO"
#  - duplicate the header row into a new row 2 with the same text, styled as
#    "not bold, bordered, centered/top" (a sub-header row)
#  - fill row 3 with sequential numbers 1..12 (bold without border on columns
#    C, F, I, L; regular without border elsewhere)
#  - update the active selection / scroll position to match the new view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L")

# 1) Rename the shared string used by L1 first (while it is the sole reference)
#    so the shared-string table keeps the same unique-string count; the new
#    row 2 below will then reuse the same (renamed) string for column L.
$ws.Range("L1").Value = "CATEGORIA_PRODUTO"

# Harmonize the header row's formatting: G1/I1/K1/L1 used a font that only
# differed from the rest of the header by an unused technical attribute
# (charset). Re-apply the same bold/border/alignment so every header cell
# shares a single, consistent style (matching the rest of row 1).
$headerFixCells = @("G1", "I1", "K1", "L1")
foreach ($c in $headerFixCells) {
    $cell = $ws.Range($c)
    $cell.Font.Bold = $true
    $cell.Font.Name = "Calibri"
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# 2) Populate row 2 with the header text repeated.
$headerValues = @("CODPRO","DESCRICAO","LINK_FORNECEDOR","DESCRICAO_TITULO","IMAGENS_BAIXADAS","CONTEUDO_DA_EMBALAGEM","CONTEUDO_DA_EMBALAGEM","DETALHES_TECNICOS","DETALHES_TECNICOS","CERTIFICADOS","CERTIFICADOS_HTML","CATEGORIA_PRODUTO")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $headerValues[$i]
}

# Style row 2 as a whole: not bold, bordered box, centered/top aligned
# (same border/alignment as row 1, but without bold).
$row2 = $ws.Range("A2:L2")
$row2.Font.Bold = $false
$row2.Borders.LineStyle = 1
$row2.HorizontalAlignment = -4108
$row2.VerticalAlignment = -4160

# 3) Populate row 3 with sequential numbers 1..12.
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $i + 1
}

# Columns C, F, I, L are bold with no border; the rest are regular, no border.
$row3 = $ws.Range("A3:L3")
$row3.Font.Bold = $false
$row3.Borders.LineStyle = 0

$boldCols3 = @("C", "F", "I", "L")
foreach ($col in $boldCols3) {
    $cell = $ws.Range($col + "3")
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 0
}

# 4) Update the selection / scroll position to match the new view.
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N2").Select() | Out-Null
